# Update crypto price/volume data per Jan 6 2024 07:20 UTC GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.796.45'
$ws.Range('E2').Value = '  -0.39%  '

$ws.Range('D3').Value = '2.223.97'
$ws.Range('E3').Value = '  -1.80%  '

$ws.Range('D4').Value = '''1.01'
$ws.Range('E4').Value = '  +0.41%  '

$ws.Range('D5').Value = '''301.14'
$ws.Range('E5').Value = '  -5.55%  '

$ws.Range('D6').Value = '''92.69'
$ws.Range('E6').Value = '  -8.91%  '

$ws.Range('E7').Value = '  -2.15%  '

$ws.Range('E8').Value = '  +0.20%  '

$ws.Range('D9').Value = '''0.511'
$ws.Range('E9').Value = '  -8.05%  '

$ws.Range('D10').Value = '''33.81'
$ws.Range('E10').Value = '  -9.13%  '

$ws.Range('D11').Value = '''0.0792'
$ws.Range('E11').Value = '  -4.83%  '

$ws.Range('D12').Value = '''7.01'
$ws.Range('E12').Value = '  -8.56%  '

$ws.Range('E13').Value = '  -3.42%  '

$ws.Range('D14').Value = '2.563.32'
$ws.Range('E14').Value = '  -1.70%  '

$ws.Range('D15').Value = '2.258.07'
$ws.Range('E15').Value = '  -0.51%  '

$ws.Range('D16').Value = '''0.802'
$ws.Range('E16').Value = '  -7.22%  '

$ws.Range('D17').Value = '''13.28'
$ws.Range('E17').Value = '  -8.43%  '

$ws.Range('D18').Value = '43.560.41'
$ws.Range('E18').Value = '  -0.72%  '

$ws.Range('D19').Value = '0.0₃0939'
$ws.Range('E19').Value = '  -4.64%  '

$ws.Range('D20').Value = '''11.88'
$ws.Range('E20').Value = '  -10.89%  '

$ws.Range('D21').Value = '''6.05'
$ws.Range('E21').Value = '  -7.85%  '

$ws.Range('D22').Value = '''63.93'
$ws.Range('E22').Value = '  -2.72%  '

$ws.Range('D23').Value = '''233.24'
$ws.Range('E23').Value = '  -0.95%  '

$ws.Range('E24').Value = '  -8.62%  '

$ws.Range('E25').Value = '  -0.12%  '

$ws.Range('D26').Value = '''1.89'
$ws.Range('E26').Value = '  -10.41%  '

$ws.Range('D27').Value = '''9.64'
$ws.Range('E27').Value = '  -4.91%  '

$ws.Range('D28').Value = '''2.11'
$ws.Range('E28').Value = '  -2.66%  '

$ws.Range('D29').Value = '''35.53'
$ws.Range('E29').Value = '  -4.85%  '

$ws.Range('D30').Value = '''5.79'
$ws.Range('E30').Value = '  -7.28%  '

$ws.Range('D31').Value = '''19.58'
$ws.Range('E31').Value = '  -3.41%  '

$ws.Range('D32').Value = '''150.53'
$ws.Range('E32').Value = '  -4.63%  '

$ws.Range('D33').Value = '''0.0788'
$ws.Range('E33').Value = '  -7.73%  '

$ws.Range('D34').Value = '''3.22'
$ws.Range('E34').Value = '  +4.30%  '

$ws.Range('E35').Value = '  -4.56%  '

$ws.Range('E36').Value = '  -2.50%  '

$ws.Range('D37').Value = '''0.105'
$ws.Range('E37').Value = '  -9.37%  '

$ws.Range('E38').Value = '  -11.95%  '

$ws.Range('D39').Value = '''14.36'
$ws.Range('E39').Value = '  -10.89%  '

$ws.Range('D40').Value = '''3.72'
$ws.Range('E40').Value = '  -11.57%  '

$ws.Range('B41').Value = 'NEARProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D41').Value = '''3.21'
$ws.Range('E41').Value = '  -13.56%  '

$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').Value = '''0.0292'
$ws.Range('E42').Value = '  -7.76%  '

$ws.Range('E43').Value = '  +0.24%  '

$ws.Range('D44').Value = '1.716.85'
$ws.Range('E44').Value = '  -4.45%  '

$ws.Range('D45').Value = '''82.28'
$ws.Range('E45').Value = '  +0.00%  '

$ws.Range('D46').Value = '''4.86'
$ws.Range('E46').Value = '  -6.81%  '

$ws.Range('D47').Value = '''0.181'
$ws.Range('E47').Value = '  -8.84%  '

$ws.Range('D48').Value = '''97.72'
$ws.Range('E48').Value = '  -6.50%  '

$ws.Range('D49').Value = '''7.95'
$ws.Range('E49').Value = '  -5.09%  '

$ws.Range('B50').Value = 'ordi'
$ws.Range('C50').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D50').Value = '''66.85'
$ws.Range('E50').Value = '  -11.77%  '

$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').Value = '''52.48'
$ws.Range('E51').Value = '  -10.21%  '
